$d = $word.ActiveDocument

# BreakAndAfterDoing
$find_BreakAndAfterDoing = "own novel solution ideas." + [char]11 + [char]11 + "After the literature review in the disciplines, as mentioned above,"
$repl_BreakAndAfterDoing = "own novel solution ideas." + [char]11 + " After doing the literature review in the disciplines mentioned above,"
$ok_BreakAndAfterDoing = $d.Content.Find.Execute($find_BreakAndAfterDoing, $false, $false, $false, $false, $false, $true, 1, $false, $repl_BreakAndAfterDoing, 2)
Write-Output ("BreakAndAfterDoing=" + $ok_BreakAndAfterDoing)

# VenuesGeocoded
$find_VenuesGeocoded = "design lessons of extensibility of columns for example, " + [char]0x2019 + "venues were geocoded to allow spatial graphs" + [char]0x2019 + " could be related as an example dates in bug reports to some standard format for all tools used and shown on a unified interface."
$repl_VenuesGeocoded = "design lessons of extensibility of columns. For example, 'venues were geocoded to allow spatial graphs' could be related as dates in bug reports to some standard format. It is done for all tools used and then shown on a unified interface."
$ok_VenuesGeocoded = $d.Content.Find.Execute($find_VenuesGeocoded, $false, $false, $false, $false, $false, $true, 1, $false, $repl_VenuesGeocoded, 2)
Write-Output ("VenuesGeocoded=" + $ok_VenuesGeocoded)

# ToySearch
$find_ToySearch = "For example, a database with toys is searched linearly for a given query it takes more time than a modified query, i.e., let us say a toy in red colour and horse is type, then search is simplified by looking at two parameters, i.e., colour and type."
$repl_ToySearch = "For example, it linearly searches a database with toys.  It takes more time than a modified query like  'a toy in red colour and horse types'.  By looking at two parameters, i.e., colour and type simplifies the search."
$ok_ToySearch = $d.Content.Find.Execute($find_ToySearch, $false, $false, $false, $false, $false, $true, 1, $false, $repl_ToySearch, 2)
Write-Output ("ToySearch=" + $ok_ToySearch)

# RefactoringToolsBuilt
$find_RefactoringToolsBuilt = "Next in the area of " + [char]0x2019 + "Refactoring tools" + [char]0x2019 + ", Dustinca talks about how the Refactoring tools are to be built"
$repl_RefactoringToolsBuilt = "Next in the area of " + [char]0x2019 + "Refactoring tools" + [char]0x2019 + ", Dustinca talks about how these tools are to be built"
$ok_RefactoringToolsBuilt = $d.Content.Find.Execute($find_RefactoringToolsBuilt, $false, $false, $false, $false, $false, $true, 1, $false, $repl_RefactoringToolsBuilt, 2)
Write-Output ("RefactoringToolsBuilt=" + $ok_RefactoringToolsBuilt)

# RefactoringToolsComma
$find_RefactoringToolsComma = "for software refactoring tools and this could perhaps give some basic guidelines"
$repl_RefactoringToolsComma = "for software refactoring tools, and this could perhaps give some basic guidelines"
$ok_RefactoringToolsComma = $d.Content.Find.Execute($find_RefactoringToolsComma, $false, $false, $false, $false, $false, $true, 1, $false, $repl_RefactoringToolsComma, 2)
Write-Output ("RefactoringToolsComma=" + $ok_RefactoringToolsComma)

# IssueTrackerOverload
$find_IssueTrackerOverload = "It is found out in their research paper that there is a too much of information they receive which confuses the developer in how to react, example: the developer receive a high number of bugs reported via email and this leads to a situation where the developer ignore the email."
$repl_IssueTrackerOverload = "It is found out in their research paper that there is too much of information they receive. It confuses the developer in how to react, for example, the developer receives a high number of bugs reported via email, and this leads to a situation where the developer ignore the email."
$ok_IssueTrackerOverload = $d.Content.Find.Execute($find_IssueTrackerOverload, $false, $false, $false, $false, $false, $true, 1, $false, $repl_IssueTrackerOverload, 2)
Write-Output ("IssueTrackerOverload=" + $ok_IssueTrackerOverload)

# ExcitingIdeasComma
$find_ExcitingIdeasComma = "They found some exciting solution ideas such as having a private dashboard"
$repl_ExcitingIdeasComma = "They found some exciting solution ideas, such as having a private dashboard"
$ok_ExcitingIdeasComma = $d.Content.Find.Execute($find_ExcitingIdeasComma, $false, $false, $false, $false, $false, $true, 1, $false, $repl_ExcitingIdeasComma, 2)
Write-Output ("ExcitingIdeasComma=" + $ok_ExcitingIdeasComma)

# ExpressivenessComma
$find_ExpressivenessComma = "Expressiveness is one other mentioned in their paper which says an example, severity or priority are vague terms"
$repl_ExpressivenessComma = "Expressiveness is one other mentioned in their paper, which says an example, severity or priority are vague terms"
$ok_ExpressivenessComma = $d.Content.Find.Execute($find_ExpressivenessComma, $false, $false, $false, $false, $false, $true, 1, $false, $repl_ExpressivenessComma, 2)
Write-Output ("ExpressivenessComma=" + $ok_ExpressivenessComma)

# StackOverflowSentence
$find_StackOverflowSentence = "Next in " + [char]0x2019 + "Stack Overflow" + [char]0x2019 + ", in a research paper by Wang et al.  it is found there are 10934198 questions on a " + [char]0x2019 + "User Interface" + [char]0x2019 + " topic for example. It is quite challenging to go through such a high volume database, but the Stack Overflow team has a friendly user interface,"
$repl_StackOverflowSentence = "Next in 'Stack Overflow', in a research paper by Wang et al.  it is found there are 10934198 questions on a 'User Interface' topic, for example. It is quite challenging to go through such a high volume database. However, the Stack Overflow team has a friendly user interface,"
$ok_StackOverflowSentence = $d.Content.Find.Execute($find_StackOverflowSentence, $false, $false, $false, $false, $false, $true, 1, $false, $repl_StackOverflowSentence, 2)
Write-Output ("StackOverflowSentence=" + $ok_StackOverflowSentence)
